# feat: add 2022-Q3 data
#
# Before: sheets = [ "总计", "2022-Q1" ]
# After:  sheets = [ "总计", "2022-Q3", "2022-Q1" ]
#   - "总计" gets a new row (2022-Q3 totals), old 2022-Q1 row moves down one row.
#   - The original "2022-Q1" detail sheet is cloned unchanged into a new sheet
#     which keeps the name "2022-Q1".
#   - The original "2022-Q1" detail sheet object is repurposed (renamed +
#     its data replaced) to become the new "2022-Q3" detail sheet.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)   # "总计"
$q1    = $wb.Worksheets.Item(2)   # "2022-Q1" (existing)

# ---------------------------------------------------------------------
# 1. Clone the existing "2022-Q1" sheet so its data survives unchanged.
#    The clone is placed right after it and will keep the "2022-Q1" name.
# ---------------------------------------------------------------------
$q1.Copy($null, $q1)
$q1Clone = $wb.Worksheets.Item(3)

# Rename the original first (frees up the "2022-Q1" name), then rename
# the clone back to "2022-Q1".
$q1.Name = "2022-Q3"
$q1Clone.Name = "2022-Q1"

# ---------------------------------------------------------------------
# 2. Update "总计": insert the 2022-Q3 totals above the existing 2022-Q1
#    row (which moves from row 2 to row 3).
# ---------------------------------------------------------------------

# Give the (future) row 3 the same style the index cell in row 2 has,
# then move the old 2022-Q1 values down into row 3.
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)   # xlPasteFormats

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q1"
$total.Cells.Item(3, 3).Value = 1
$total.Cells.Item(3, 4).Value = 0.17

# Overwrite row 2 in place with the new 2022-Q3 totals (keeps its style).
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 4
$total.Cells.Item(2, 4).Value = 0.07

# ---------------------------------------------------------------------
# 3. Rebuild the "2022-Q3" detail sheet (currently still holding the old
#    "2022-Q1" fund data) with the new fund table.
# ---------------------------------------------------------------------
$q1.Cells.Clear()

# Header formatting/border/alignment copied from "总计"'s header style.
$total.Range("B1:D1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

# Index-column formatting copied from "总计"'s styled numeric cell.
$total.Range("A2").Copy()
$q1.Range("A2:A5").PasteSpecial(-4122)   # xlPasteFormats

# Header row text.
$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

# Force text storage for columns that must preserve exact string content
# (fund codes with leading zeros, and numeric-looking ratio/size figures).
$q1.Range("B2:B5").NumberFormat = "@"
$q1.Range("D2:G5").NumberFormat = "@"

# Row 2
$q1.Cells.Item(2, 1).Value = 0
$q1.Cells.Item(2, 2).Value = "011404"
$q1.Cells.Item(2, 3).Value = "融通鑫新成长混合C"
$q1.Cells.Item(2, 4).Value = "1.75"
$q1.Cells.Item(2, 5).Value = "94.07"
$q1.Cells.Item(2, 6).Value = "2.89"
$q1.Cells.Item(2, 7).Value = "0.0506"
$q1.Cells.Item(2, 8).Value = 10

# Row 3
$q1.Cells.Item(3, 1).Value = 1
$q1.Cells.Item(3, 2).Value = "011403"
$q1.Cells.Item(3, 3).Value = "融通鑫新成长混合A"
$q1.Cells.Item(3, 4).Value = "0.39"
$q1.Cells.Item(3, 5).Value = "94.07"
$q1.Cells.Item(3, 6).Value = "2.89"
$q1.Cells.Item(3, 7).Value = "0.0113"
$q1.Cells.Item(3, 8).Value = 10

# Row 4
$q1.Cells.Item(4, 1).Value = 2
$q1.Cells.Item(4, 2).Value = "013869"
$q1.Cells.Item(4, 3).Value = "创金合信物联网主题股票A"
$q1.Cells.Item(4, 4).Value = "0.19"
$q1.Cells.Item(4, 5).Value = "83.68"
$q1.Cells.Item(4, 6).Value = "3.20"
$q1.Cells.Item(4, 7).Value = "0.0061"
$q1.Cells.Item(4, 8).Value = 9

# Row 5
$q1.Cells.Item(5, 1).Value = 3
$q1.Cells.Item(5, 2).Value = "013870"
$q1.Cells.Item(5, 3).Value = "创金合信物联网主题股票C"
$q1.Cells.Item(5, 4).Value = "0.13"
$q1.Cells.Item(5, 5).Value = "83.68"
$q1.Cells.Item(5, 6).Value = "3.20"
$q1.Cells.Item(5, 7).Value = "0.0042"
$q1.Cells.Item(5, 8).Value = 9

# Match the page margins used by the newly generated sheets (0.75"/1"/0.5").
$q1.PageSetup.LeftMargin = 0.75 * 72
$q1.PageSetup.RightMargin = 0.75 * 72
$q1.PageSetup.TopMargin = 1 * 72
$q1.PageSetup.BottomMargin = 1 * 72
$q1.PageSetup.HeaderMargin = 0.5 * 72
$q1.PageSetup.FooterMargin = 0.5 * 72
